# Update metrics values for rows 2-26 (columns B..Q) with new computed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New values for each metric column (same values applied to every data row)
$values = @{
    "B" = 0.445980542904741
    "C" = -0.05067851974480586
    "D" = 0.6167507148584341
    "E" = 0.06596443554798825
    "F" = 0.4591084616579356
    "G" = 0.3288898821261522
    "H" = 0.6237281562693088
    "I" = 0.3771784063594097
    "J" = 0.4532621527071558
    "K" = 0.4152202795332827
    "L" = 0.279043433096329
    "M" = 0.5734892170966707
    "N" = 0.05025235926527027
    "O" = 0.5979038403548141
    "P" = 22.22406457805329
    "Q" = 34.41282282673529
}

$columns = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")

for ($row = 2; $row -le 26; $row++) {
    foreach ($col in $columns) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
